# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Sat Oct  5 00:00:12 UTC 2024 with GitHub Actions".
#
# The sheet stores Price (column D) and Volume(1h) (column E) as plain text
# (inline strings), even though many Price values look numeric. To avoid Excel
# silently re-interpreting strings such as "5.40" or "1.00" as numbers (which
# would drop the trailing zero / formatting), column D cells are explicitly
# formatted as Text ("@") before the new value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.101.80'
$ws.Range("E2").Value = '  +2.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.416.98'
$ws.Range("E3").Value = '  +2.85%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.33'
$ws.Range("E5").Value = '  +1.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.22'
$ws.Range("E6").Value = '  +4.64%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  +2.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.416.82'
$ws.Range("E9").Value = '  +2.87%  '

$ws.Range("E10").Value = '  +4.02%  '

$ws.Range("E11").Value = '  -0.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.40'
$ws.Range("E12").Value = '  +1.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").Value = '  +1.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.25'
$ws.Range("E14").Value = '  +6.48%  '

$ws.Range("E15").Value = '  +8.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.856.25'
$ws.Range("E16").Value = '  +2.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.969.33'
$ws.Range("E17").Value = '  +1.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.418.86'
$ws.Range("E18").Value = '  +2.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.08'
$ws.Range("E19").Value = '  +3.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.20'
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.86'
$ws.Range("E21").Value = '  +1.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.72'
$ws.Range("E22").Value = '  +2.47%  '

$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("E24").Value = '  +6.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.98'
$ws.Range("E25").Value = '  +2.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.10'
$ws.Range("E26").Value = '  +6.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '585.31'
$ws.Range("E27").Value = '  +17.75%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.536.85'
$ws.Range("E28").Value = '  +2.92%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0940'
$ws.Range("E30").Value = '  +8.83%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.36'
$ws.Range("E31").Value = '  +4.29%  '

$ws.Range("E32").Value = '  +6.99%  '

$ws.Range("E33").Value = '  +1.21%  '

$ws.Range("E34").Value = '  +3.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("E35").Value = '  +3.04%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.68'
$ws.Range("E37").Value = '  +8.18%  '

$ws.Range("E38").Value = '  +4.10%  '

$ws.Range("E39").Value = '  +2.21%  '

$ws.Range("E40").Value = '  +2.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.78'
$ws.Range("E41").Value = '  +1.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '148.81'
$ws.Range("E42").Value = '  +2.85%  '

$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.79'
$ws.Range("E44").Value = '  +2.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '151.25'
$ws.Range("E45").Value = '  +5.49%  '

$ws.Range("E46").Value = '  +12.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.65'
$ws.Range("E47").Value = '  +2.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0547'
$ws.Range("E48").Value = '  +5.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.39'
$ws.Range("E49").Value = '  +7.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.589'
$ws.Range("E50").Value = '  +3.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0918'
$ws.Range("E51").Value = '  +1.62%  '

